$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 4996708.5
$ws.Cells.Item(33, 9).Value = 7687131
$ws.Cells.Item(33, 10).Value = 210.28572
$ws.Cells.Item(33, 11).Value = 7687131
$ws.Cells.Item(33, 12).Value = 210.28572
$ws.Cells.Item(33, 13).Value = -7686902
$ws.Cells.Item(33, 14).Value = -668.28572
$ws.Cells.Item(40, 8).Value = 1421.1818
$ws.Cells.Item(40, 9).Value = 1536
$ws.Cells.Item(40, 10).Value = 1355.5714
$ws.Cells.Item(40, 11).Value = 1536
$ws.Cells.Item(40, 12).Value = 1355.5714
$ws.Cells.Item(40, 13).Value = -1361
$ws.Cells.Item(40, 14).Value = -1705.5714
$ws.Cells.Item(51, 8).Value = 2640
$ws.Cells.Item(51, 10).Value = 4000
$ws.Cells.Item(51, 12).Value = 4000
$ws.Cells.Item(51, 14).Value = -4968
$ws.Cells.Item(107, 8).Value = 531.7778
$ws.Cells.Item(107, 10).Value = 421.66666
$ws.Cells.Item(107, 12).Value = 421.66666
$ws.Cells.Item(107, 14).Value = -4261.66666
$ws.Cells.Item(111, 8).Value = 2406.9285
$ws.Cells.Item(111, 9).Value = 2457
$ws.Cells.Item(111, 10).Value = 2356.8572
$ws.Cells.Item(111, 11).Value = 7371
$ws.Cells.Item(111, 12).Value = 7070.571599999999
$ws.Cells.Item(111, 13).Value = -4304
$ws.Cells.Item(111, 14).Value = -13204.5716
$ws.Cells.Item(116, 8).Value = 3155.0527
$ws.Cells.Item(116, 9).Value = 3620
$ws.Cells.Item(116, 10).Value = 2638.4443
$ws.Cells.Item(116, 11).Value = 3620
$ws.Cells.Item(116, 12).Value = 2638.4443
$ws.Cells.Item(116, 13).Value = -178
$ws.Cells.Item(116, 14).Value = -9522.444299999999
$ws.Cells.Item(132, 8).Value = 1306.3889
$ws.Cells.Item(132, 9).Value = 1250.5625
$ws.Cells.Item(132, 11).Value = 3751.6875
$ws.Cells.Item(132, 13).Value = -1221.6875
$ws.Cells.Item(137, 8).Value = 982696.7
$ws.Cells.Item(137, 9).Value = 2922.65
$ws.Cells.Item(137, 10).Value = 2382374
$ws.Cells.Item(137, 11).Value = 8767.950000000001
$ws.Cells.Item(137, 12).Value = 7147122
$ws.Cells.Item(137, 13).Value = -6217.950000000001
$ws.Cells.Item(137, 14).Value = -7152222

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 15735.333
$ws.Cells.Item(31, 9).Value = 15735.333
$ws.Cells.Item(31, 11).Value = 15735.333
$ws.Cells.Item(31, 13).Value = -15441.333
$ws.Cells.Item(32, 8).Value = 18533.613
$ws.Cells.Item(32, 9).Value = 23309.979
$ws.Cells.Item(32, 10).Value = 3567.6667
$ws.Cells.Item(32, 11).Value = 23309.979
$ws.Cells.Item(32, 12).Value = 3567.6667
$ws.Cells.Item(32, 13).Value = -23022.979
$ws.Cells.Item(32, 14).Value = -4141.6667
$ws.Cells.Item(45, 8).Value = 1801.3784
$ws.Cells.Item(45, 9).Value = 1817.48
$ws.Cells.Item(45, 11).Value = 1817.48
$ws.Cells.Item(45, 13).Value = -1440.48
$ws.Cells.Item(62, 8).Value = 38249
$ws.Cells.Item(62, 10).Value = 38249
$ws.Cells.Item(62, 12).Value = 38249
$ws.Cells.Item(62, 14).Value = -39497
$ws.Cells.Item(65, 8).Value = 38249
$ws.Cells.Item(65, 10).Value = 38249
$ws.Cells.Item(65, 12).Value = 114747
$ws.Cells.Item(65, 14).Value = -120987
$ws.Cells.Item(110, 8).Value = 1932.7894
$ws.Cells.Item(110, 9).Value = 1793.5714
$ws.Cells.Item(110, 11).Value = 1793.5714
$ws.Cells.Item(110, 13).Value = 251.4286
$ws.Cells.Item(122, 8).Value = 41668000
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 13).Value = -3550

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1467.3572
$ws.Cells.Item(99, 9).Value = 977.2778
$ws.Cells.Item(99, 10).Value = 2349.5
$ws.Cells.Item(99, 11).Value = 977.2778
$ws.Cells.Item(99, 12).Value = 2349.5
$ws.Cells.Item(99, 13).Value = 520.7222
$ws.Cells.Item(99, 14).Value = -5345.5
$ws.Cells.Item(134, 8).Value = 148003.72
$ws.Cells.Item(134, 9).Value = 5202.4
$ws.Cells.Item(134, 10).Value = 505007
$ws.Cells.Item(134, 11).Value = 15607.2
$ws.Cells.Item(134, 12).Value = 1515021
$ws.Cells.Item(134, 13).Value = -13072.2
$ws.Cells.Item(134, 14).Value = -1520091

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2760092.8
$ws.Cells.Item(58, 9).Value = 3789630
$ws.Cells.Item(58, 10).Value = 14659.333
$ws.Cells.Item(58, 11).Value = 3789630
$ws.Cells.Item(58, 12).Value = 14659.333
$ws.Cells.Item(58, 13).Value = -3789427
$ws.Cells.Item(58, 14).Value = -15065.333
$ws.Cells.Item(76, 8).Value = 9192.714
$ws.Cells.Item(76, 9).Value = 9192.714
$ws.Cells.Item(76, 11).Value = 9192.714
$ws.Cells.Item(76, 13).Value = -8877.714
$ws.Cells.Item(79, 8).Value = 9192.714
$ws.Cells.Item(79, 9).Value = 9192.714
$ws.Cells.Item(79, 11).Value = 9192.714
$ws.Cells.Item(79, 13).Value = -8100.714
$ws.Cells.Item(99, 8).Value = 1283.4117
$ws.Cells.Item(99, 9).Value = 1307.6
$ws.Cells.Item(99, 10).Value = 1102
$ws.Cells.Item(99, 11).Value = 1307.6
$ws.Cells.Item(99, 12).Value = 1102
$ws.Cells.Item(99, 13).Value = 190.4000000000001
$ws.Cells.Item(99, 14).Value = -4098
$ws.Cells.Item(105, 8).Value = 609.1667
$ws.Cells.Item(105, 9).Value = 621.4706
$ws.Cells.Item(105, 10).Value = 400
$ws.Cells.Item(105, 11).Value = 621.4706
$ws.Cells.Item(105, 12).Value = 400
$ws.Cells.Item(105, 13).Value = 1125.5294
$ws.Cells.Item(105, 14).Value = -3894
$ws.Cells.Item(122, 8).Value = 5979.3228
$ws.Cells.Item(122, 9).Value = 2502.8262
$ws.Cells.Item(122, 10).Value = 15974.25
$ws.Cells.Item(122, 11).Value = 7508.4786
$ws.Cells.Item(122, 12).Value = 47922.75
$ws.Cells.Item(122, 13).Value = -5058.4786
$ws.Cells.Item(122, 14).Value = -52822.75
$ws.Cells.Item(126, 8).Value = 1283.4117
$ws.Cells.Item(126, 9).Value = 1307.6
$ws.Cells.Item(126, 10).Value = 1102
$ws.Cells.Item(126, 11).Value = 3922.8
$ws.Cells.Item(126, 12).Value = 3306
$ws.Cells.Item(126, 13).Value = -1452.8
$ws.Cells.Item(126, 14).Value = -8246
$ws.Cells.Item(132, 8).Value = 4535.2173
$ws.Cells.Item(132, 9).Value = 4406.278
$ws.Cells.Item(132, 10).Value = 4999.4
$ws.Cells.Item(132, 11).Value = 13218.834
$ws.Cells.Item(132, 12).Value = 14998.2
$ws.Cells.Item(132, 13).Value = -10688.834
$ws.Cells.Item(132, 14).Value = -20058.2
$ws.Cells.Item(134, 8).Value = 3091.9583
$ws.Cells.Item(134, 9).Value = 2632.1667
$ws.Cells.Item(134, 10).Value = 4471.3335
$ws.Cells.Item(134, 11).Value = 7896.500100000001
$ws.Cells.Item(134, 12).Value = 13414.0005
$ws.Cells.Item(134, 13).Value = -5361.500100000001
$ws.Cells.Item(134, 14).Value = -18484.0005
$ws.Cells.Item(136, 8).Value = 2760092.8
$ws.Cells.Item(136, 9).Value = 3789630
$ws.Cells.Item(136, 10).Value = 14659.333
$ws.Cells.Item(136, 11).Value = 11368890
$ws.Cells.Item(136, 12).Value = 43977.999
$ws.Cells.Item(136, 13).Value = -11366340
$ws.Cells.Item(136, 14).Value = -49077.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 129575.125
$ws.Cells.Item(56, 9).Value = 129575.125
$ws.Cells.Item(56, 11).Value = 129575.125
$ws.Cells.Item(56, 13).Value = -129045.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2666.7827
$ws.Cells.Item(126, 9).Value = 1778
$ws.Cells.Item(126, 11).Value = 5334
$ws.Cells.Item(126, 13).Value = -2864

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1500
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -705
$ws.Cells.Item(22, 14).Value = -2590
$ws.Cells.Item(27, 8).Value = 1500
$ws.Cells.Item(27, 9).Value = 1000
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = -893
$ws.Cells.Item(27, 14).Value = -2214
$ws.Cells.Item(46, 8).Value = 772.38464
$ws.Cells.Item(46, 9).Value = 615.6667
$ws.Cells.Item(46, 10).Value = 1125
$ws.Cells.Item(46, 11).Value = 615.6667
$ws.Cells.Item(46, 12).Value = 1125
$ws.Cells.Item(46, 13).Value = -427.6667
$ws.Cells.Item(46, 14).Value = -1501
$ws.Cells.Item(55, 8).Value = 500131.12
$ws.Cells.Item(55, 9).Value = 666766.5
$ws.Cells.Item(55, 10).Value = 225
$ws.Cells.Item(55, 11).Value = 666766.5
$ws.Cells.Item(55, 12).Value = 225
$ws.Cells.Item(55, 13).Value = -666593.5
$ws.Cells.Item(55, 14).Value = -571
$ws.Cells.Item(133, 8).Value = 46243
$ws.Cells.Item(133, 10).Value = 46243
$ws.Cells.Item(133, 12).Value = 46243
$ws.Cells.Item(133, 14).Value = -51303

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 22071
$ws.Cells.Item(68, 10).Value = 22071
$ws.Cells.Item(68, 12).Value = 22071
$ws.Cells.Item(68, 14).Value = -23693
$ws.Cells.Item(71, 8).Value = 22071
$ws.Cells.Item(71, 10).Value = 22071
$ws.Cells.Item(71, 12).Value = 66213
$ws.Cells.Item(71, 14).Value = -74325
$ws.Cells.Item(107, 8).Value = 1633.6666
$ws.Cells.Item(107, 9).Value = 1184
$ws.Cells.Item(107, 10).Value = 1993.4
$ws.Cells.Item(107, 11).Value = 3552
$ws.Cells.Item(107, 12).Value = 5980.200000000001
$ws.Cells.Item(107, 13).Value = -1632
$ws.Cells.Item(107, 14).Value = -9820.200000000001
$ws.Cells.Item(132, 8).Value = 3792.5293
$ws.Cells.Item(132, 9).Value = 3748.1428
$ws.Cells.Item(132, 10).Value = 3999.6667
$ws.Cells.Item(132, 11).Value = 11244.4284
$ws.Cells.Item(132, 12).Value = 11999.0001
$ws.Cells.Item(132, 13).Value = -8714.428400000001
$ws.Cells.Item(132, 14).Value = -17059.0001
